# Insert a new data record at row 278 (shifting the existing rows 278..366
# down to 279..367) and populate it with the new Coliflor / Macroferia
# Regional de Talca observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 278 (and everything below it) down by one row.
$ws.Rows.Item(278).Insert()

# Fill in the newly inserted row 278 with the new record's values. The
# columns that don't vary row-to-row in this sheet (A, B, C, E, F, G, H, I,
# N, O, Q, R) keep the same values as their neighbouring rows.
$ws.Cells.Item(278, 1).Value2 = 5
$ws.Cells.Item(278, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(278, 3).Value2 = "Maule"
$ws.Cells.Item(278, 4).Value2 = 44876
$ws.Cells.Item(278, 5).Value2 = 7
$ws.Cells.Item(278, 6).Value2 = 100112008
$ws.Cells.Item(278, 7).Value2 = "Coliflor"
$ws.Cells.Item(278, 8).Value2 = "Sin especificar"
$ws.Cells.Item(278, 9).Value2 = "Primera"
$ws.Cells.Item(278, 10).Value2 = 3000
$ws.Cells.Item(278, 11).Value2 = 900
$ws.Cells.Item(278, 12).Value2 = 900
$ws.Cells.Item(278, 13).Value2 = 900
$ws.Cells.Item(278, 14).Value2 = "`$/unidad"
$ws.Cells.Item(278, 15).Value2 = "Región del Maule"
$ws.Cells.Item(278, 16).Value2 = 900
$ws.Cells.Item(278, 17).Value2 = 1
$ws.Cells.Item(278, 18).Value2 = "Hortaliza"

# Keep the date column's custom date/time number format on the new row,
# matching the style used by the rest of the "Fecha" column.
$ws.Cells.Item(278, 4).NumberFormat = $ws.Cells.Item(279, 4).NumberFormat
